# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values computed for each outing row (rows 2-21)
$s_vals = @{
    2  = 3
    3  = 3
    4  = 5
    5  = 4
    6  = 5
    7  = 4
    8  = 3
    9  = 4
    10 = 2
    11 = 6
    12 = 3
    13 = 9
    14 = 4
    15 = 0
    16 = 6
    17 = 1
    18 = 1
    19 = 3
    20 = 0
    21 = 3
}

foreach ($row in $s_vals.Keys) {
    $ws.Cells.Item($row, 7).Value = $s_vals[$row]
}
